$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.08053555155723416
$ws.Range("C2").Value = 0.2607125751647556
$ws.Range("B3").Value = 0.1023820116553337
$ws.Range("C3").Value = -0.2643531791085623
$ws.Range("B4").Value = 0.2049684770545157
$ws.Range("C4").Value = 0.2225355073906286
$ws.Range("B5").Value = 0.2169909156117429
$ws.Range("C5").Value = -0.1084341318948621
$ws.Range("B6").Value = 0.1444011345976094
$ws.Range("C6").Value = 0.2319578760726343
$ws.Range("B7").Value = 0.3332653301064796
$ws.Range("C7").Value = 0.39392106858569
$ws.Range("B8").Value = 0.318050236718135
$ws.Range("C8").Value = -0.4414046376317643
$ws.Range("B9").Value = 0.1051321814098121
$ws.Range("C9").Value = -0.1550931849627726
$ws.Range("B10").Value = 0.1189829256514251
$ws.Range("C10").Value = 0.2059546282531164
$ws.Range("B11").Value = 0.4547933661645938
$ws.Range("C11").Value = -0.0505843682346436
$ws.Range("B12").Value = 0.4151417231151884
$ws.Range("C12").Value = -0.1745899015960837
$ws.Range("B13").Value = 0.350448853305261
$ws.Range("C13").Value = 0.1909020037940034
$ws.Range("B14").Value = 0.2971744291466285
$ws.Range("C14").Value = 0.04165925157841546
$ws.Range("B15").Value = 0.1761520357767008
$ws.Range("C15").Value = -0.1552805762572122
$ws.Range("B16").Value = 0.102652471975735
$ws.Range("C16").Value = 0.07281359419034758
$ws.Range("B17").Value = 0.06146474928493108
$ws.Range("C17").Value = 0.4795262829501432
